$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: "Odd_CS_3-3_HT" column moved from BC to AW, shifting
#     AW..BC one column to the right. BD1 ("Odd_CS_4-4_HT") is unchanged.
$ws.Range("AW1").Value = "Odd_CS_3-3_HT"
$ws.Range("AX1").Value = "Odd_CS_0-1_HT"
$ws.Range("AY1").Value = "Odd_CS_0-2_HT"
$ws.Range("AZ1").Value = "Odd_CS_1-2_HT"
$ws.Range("BA1").Value = "Odd_CS_0-3_HT"
$ws.Range("BB1").Value = "Odd_CS_1-3_HT"
$ws.Range("BC1").Value = "Odd_CS_2-3_HT"

# --- Row 2: replaced with a new match's data (same columns A2:BD2).
$ws.Range("A2").Value = "xpWp1ROi"
$ws.Range("C2").Value = "05:35"
$ws.Range("D2").Value = "AUSTRALIA - A-LEAGUE"
$ws.Range("E2").Value = "WS Wanderers"
$ws.Range("F2").Value = "Newcastle Jets"

$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 2.3
$ws.Range("K2").Value = 2.5
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.63
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 11
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 13
$ws.Range("AB2").Value = 19
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 101
$ws.Range("AH2").Value = 17
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 13
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 26
$ws.Range("AM2").Value = 26
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 9
$ws.Range("AP2").Value = 15
$ws.Range("AQ2").Value = 26
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 81
$ws.Range("AT2").Value = 3.75
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 301
$ws.Range("AX2").Value = 6
$ws.Range("AY2").Value = 19
$ws.Range("AZ2").Value = 21
$ws.Range("BA2").Value = 51
$ws.Range("BB2").Value = 67
$ws.Range("BC2").Value = 101
$ws.Range("BD2").Value = 151
